$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The option-pool rows (6 and 7) get their own distinct email addresses
# instead of reusing emp3 / emp4's addresses.
$ws.Range("C6").Value = "emp5@mycompany.com"
$ws.Range("C7").Value = "emp6@mycompany.com"

# Rebuild every hyperlink so the mailto targets / display text line up with
# the (new) email addresses. Deleting a single hyperlink removes the whole
# collection in this environment, so recreate the full set in order,
# preserving the original rId ordering.
$ws.Range("C6").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:emp1@mycompany.com", "", "", "emp1@mycompany.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:emp2@mycompany.com", "", "", "emp2@mycompany.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:emp3@mycompany.com", "", "", "emp3@mycompany.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:emp4@mycompany.com", "", "", "emp4@mycompany.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:emp5@mycompany.com", "", "", "emp5@mycompany.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:emp6@mycompany.com", "", "", "emp6@mycompany.com") | Out-Null

# Hyperlinks.Add stamps the generic built-in "Hyperlink" style on the cells it
# touches; restore the workbook's original custom (blue, non-underlined,
# 11pt Arial) formatting for the whole Email column.
$rng = $ws.Range("C2:C7")
$rng.Font.Name = "Arial"
$rng.Font.Family = 1
$rng.Font.Underline = -4142
$rng.Font.Size = 11
$rng.Font.Color = 16711680

# Move the active selection to C7, as recorded in the sheet view.
$ws.Range("C7").Select() | Out-Null
